$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2962962962962963
$ws.Range("C2").Value = 0.3703703703703703
$ws.Range("P2").Value = 0.2222222222222222
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("P3").Value = 0.8888888888888888
$ws.Range("S3").Value = 0.1111111111111111
$ws.Range("P4").Value = 1
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("F6").Value = 0.03333333333333333
$ws.Range("J6").Value = 0.1
$ws.Range("Q6").Value = 0.06666666666666667
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.6666666666666666
$ws.Range("B7").Value = 0.03125
$ws.Range("F7").Value = 0.03125
$ws.Range("J7").Value = 0.125
$ws.Range("Q7").Value = 0.15625
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.59375
$ws.Range("B8").Value = 0.0660377358490566
$ws.Range("D8").Value = 0.009433962264150943
$ws.Range("F8").Value = 0.09433962264150944
$ws.Range("J8").Value = 0.09433962264150944
$ws.Range("O8").Value = 0.01886792452830189
$ws.Range("Q8").Value = 0.1792452830188679
$ws.Range("R8").Value = 0.04716981132075472
$ws.Range("S8").Value = 0.4905660377358491
$ws.Range("B9").Value = 0.06122448979591837
$ws.Range("F9").Value = 0.08163265306122448
$ws.Range("J9").Value = 0.08163265306122448
$ws.Range("O9").Value = 0.04081632653061224
$ws.Range("Q9").Value = 0.1020408163265306
$ws.Range("R9").Value = 0.06122448979591837
$ws.Range("S9").Value = 0.5714285714285714
$ws.Range("B10").Value = 0.04545454545454546
$ws.Range("D10").Value = 0.007575757575757576
$ws.Range("F10").Value = 0.08333333333333333
$ws.Range("J10").Value = 0.1136363636363636
$ws.Range("Q10").Value = 0.1742424242424243
$ws.Range("R10").Value = 0.03787878787878788
$ws.Range("S10").Value = 0.5378787878787878
$ws.Range("G11").Value = 0.2666666666666667
$ws.Range("J11").Value = 0.02222222222222222
$ws.Range("K11").Value = 0.2666666666666667
$ws.Range("L11").Value = 0.4222222222222222
$ws.Range("S11").Value = 0.02222222222222222
$ws.Range("G12").Value = 0.7894736842105263
$ws.Range("J12").Value = 0.1052631578947368
$ws.Range("S12").Value = 0.1052631578947368
$ws.Range("J13").Value = 0.1111111111111111
$ws.Range("S13").Value = 0.2222222222222222
$ws.Range("J14").Value = 1
$ws.Range("H15").Value = 0.16
$ws.Range("I15").Value = 0.16
$ws.Range("J15").Value = 0.32
$ws.Range("K15").Value = 0.08
$ws.Range("S15").Value = 0.28
$ws.Range("H16").Value = 0.3125
$ws.Range("I16").Value = 0.125
$ws.Range("J16").Value = 0.3125
$ws.Range("M16").Value = 0.0625
$ws.Range("S16").Value = 0.1875
$ws.Range("H17").Value = 0.2075471698113208
$ws.Range("I17").Value = 0.1886792452830189
$ws.Range("J17").Value = 0.2830188679245283
$ws.Range("K17").Value = 0.1132075471698113
$ws.Range("O17").Value = 0.05660377358490566
$ws.Range("S17").Value = 0.1509433962264151
$ws.Range("H18").Value = 0.2941176470588235
$ws.Range("I18").Value = 0.3529411764705883
$ws.Range("J18").Value = 0.1176470588235294
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("N18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.05882352941176471
$ws.Range("F19").Value = 0.00392156862745098
$ws.Range("H19").Value = 0.3215686274509804
$ws.Range("I19").Value = 0.1058823529411765
$ws.Range("J19").Value = 0.2549019607843137
$ws.Range("K19").Value = 0.08627450980392157
$ws.Range("M19").Value = 0.03529411764705882
$ws.Range("O19").Value = 0.04705882352941176
$ws.Range("S19").Value = 0.1450980392156863
